# GPLIM-1516: Get all tests passing!
# Update the billing/quantity figures on the single data sheet and move the
# window's selection/scroll position, matching the authoritative edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 -------------------------------------------------------------
$ws.Range("N3").Value = 0.5
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0

# --- Row 4 -------------------------------------------------------------
$ws.Range("N4").Value = 0.5
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 1

# --- Row 5 -------------------------------------------------------------
$ws.Range("N5").Value = 0.5
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 1

# --- Row 6 -------------------------------------------------------------
$ws.Range("N6").Value = 0.5
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0

# --- Window state: selection moves from N5 to N7, view scrolls right ---
$ws.Activate()
$ws.Range("N7").Select()
$excel.ActiveWindow.ScrollColumn = 19   # column S
$excel.ActiveWindow.ScrollRow = 1
